$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: "Marking" right/wrong values
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12: "Total" right value and the "X / Y" summary string
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100 / 112"
